$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 600.6957
$ws.Range("I55").Value = 697.61536
$ws.Range("J55").Value = 474.7
$ws.Range("K55").Value = 697.61536
$ws.Range("L55").Value = 474.7
$ws.Range("M55").Value = -483.61536
$ws.Range("N55").Value = -902.7
$ws.Range("H87").Value = 21363.637
$ws.Range("J87").Value = 21363.637
$ws.Range("L87").Value = 21363.637
$ws.Range("N87").Value = -23859.637
$ws.Range("H90").Value = 21363.637
$ws.Range("J90").Value = 21363.637
$ws.Range("L90").Value = 64090.91099999999
$ws.Range("N90").Value = -76570.91099999999
$ws.Range("H138").Value = 2943.32
$ws.Range("I138").Value = 1319.6666
$ws.Range("J138").Value = 4849.3477
$ws.Range("K138").Value = 3958.9998
$ws.Range("L138").Value = 14548.0431
$ws.Range("M138").Value = 1181.0002
$ws.Range("N138").Value = -24828.0431

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 13306.4
$ws.Range("I39").Value = 6633
$ws.Range("J39").Value = 40000
$ws.Range("K39").Value = 6633
$ws.Range("L39").Value = 40000
$ws.Range("M39").Value = -6113
$ws.Range("N39").Value = -41040
$ws.Range("H45").Value = 3419.625
$ws.Range("I45").Value = 2823.4167
$ws.Range("K45").Value = 2823.4167
$ws.Range("M45").Value = -2446.4167
$ws.Range("H61").Value = 2430.8572
$ws.Range("I61").Value = 1591.3
$ws.Range("K61").Value = 1591.3
$ws.Range("M61").Value = -1379.3
$ws.Range("H74").Value = 2016.2858
$ws.Range("I74").Value = 1293.4
$ws.Range("K74").Value = 1293.4
$ws.Range("M74").Value = -419.4000000000001
$ws.Range("H77").Value = 2016.2858
$ws.Range("I77").Value = 1293.4
$ws.Range("K77").Value = 6467
$ws.Range("M77").Value = -2099
$ws.Range("H97").Value = 5053151
$ws.Range("I97").Value = 1942.5264
$ws.Range("K97").Value = 1942.5264
$ws.Range("M97").Value = -1446.5264
$ws.Range("H110").Value = 1141.4
$ws.Range("I110").Value = 1181.0588
$ws.Range("K110").Value = 1181.0588
$ws.Range("M110").Value = 863.9412
$ws.Range("H122").Value = 2712.6365
$ws.Range("I122").Value = 2127.5715
$ws.Range("J122").Value = 3736.5
$ws.Range("K122").Value = 6382.7145
$ws.Range("L122").Value = 11209.5
$ws.Range("M122").Value = -3932.7145
$ws.Range("N122").Value = -16109.5
$ws.Range("H132").Value = 1348
$ws.Range("I132").Value = 1296.7333
$ws.Range("K132").Value = 3890.199900000001
$ws.Range("M132").Value = -1360.199900000001
$ws.Range("H136").Value = 2430.8572
$ws.Range("I136").Value = 1591.3
$ws.Range("K136").Value = 4773.9
$ws.Range("M136").Value = -2223.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3755.7778
$ws.Range("I86").Value = 2619
$ws.Range("J86").Value = 4424.4707
$ws.Range("K86").Value = 2619
$ws.Range("L86").Value = 4424.4707
$ws.Range("M86").Value = -1496
$ws.Range("N86").Value = -6670.4707
$ws.Range("H89").Value = 3755.7778
$ws.Range("I89").Value = 2619
$ws.Range("J89").Value = 4424.4707
$ws.Range("K89").Value = 13095
$ws.Range("L89").Value = 22122.3535
$ws.Range("M89").Value = -7479
$ws.Range("N89").Value = -33354.3535
$ws.Range("H94").Value = 15155559
$ws.Range("I94").Value = 4786.778
$ws.Range("J94").Value = 83334030
$ws.Range("K94").Value = 4786.778
$ws.Range("L94").Value = 83334030
$ws.Range("M94").Value = -4335.778
$ws.Range("N94").Value = -83334932
$ws.Range("H134").Value = 4180.485
$ws.Range("I134").Value = 3860.8696
$ws.Range("J134").Value = 4915.6
$ws.Range("K134").Value = 11582.6088
$ws.Range("L134").Value = 14746.8
$ws.Range("M134").Value = -9047.6088
$ws.Range("N134").Value = -19816.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3349.6667
$ws.Range("I31").Value = 2666.625
$ws.Range("K31").Value = 2666.625
$ws.Range("M31").Value = -2371.625
$ws.Range("H34").Value = 3349.6667
$ws.Range("I34").Value = 2666.625
$ws.Range("K34").Value = 2666.625
$ws.Range("M34").Value = -2464.625
$ws.Range("H98").Value = 32997.25
$ws.Range("I98").Value = 20000
$ws.Range("J98").Value = 37329.668
$ws.Range("K98").Value = 20000
$ws.Range("L98").Value = 37329.668
$ws.Range("M98").Value = -17754
$ws.Range("N98").Value = -41821.668
$ws.Range("H99").Value = 32061724
$ws.Range("I99").Value = 8131185.5
$ws.Range("K99").Value = 8131185.5
$ws.Range("M99").Value = -8129687.5
$ws.Range("H122").Value = 429181.88
$ws.Range("I122").Value = 681624.5600000001
$ws.Range("J122").Value = 8444.111000000001
$ws.Range("K122").Value = 2044873.68
$ws.Range("L122").Value = 25332.333
$ws.Range("M122").Value = -2042423.68
$ws.Range("N122").Value = -30232.333
$ws.Range("H126").Value = 32061724
$ws.Range("I126").Value = 8131185.5
$ws.Range("K126").Value = 24393556.5
$ws.Range("M126").Value = -24391086.5
$ws.Range("H134").Value = 3280.923
$ws.Range("I134").Value = 3111
$ws.Range("K134").Value = 9333
$ws.Range("M134").Value = -6798

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 17984.5
$ws.Range("I120").Value = 17984.5
$ws.Range("K120").Value = 53953.5
$ws.Range("M120").Value = -49115.5
$ws.Range("H132").Value = 1148
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H133").Value = 3999.5
$ws.Range("I133").Value = 3999.5
$ws.Range("K133").Value = 11998.5
$ws.Range("M133").Value = -6938.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 23891568
$ws.Range("I80").Value = 187532.67
$ws.Range("K80").Value = 187532.67
$ws.Range("M80").Value = -186534.67
$ws.Range("H83").Value = 23891568
$ws.Range("I83").Value = 187532.67
$ws.Range("K83").Value = 937663.3500000001
$ws.Range("M83").Value = -932671.3500000001
$ws.Range("H99").Value = 6720.5713
$ws.Range("I99").Value = 6720.5713
$ws.Range("K99").Value = 6720.5713
$ws.Range("M99").Value = -4474.5713
$ws.Range("H102").Value = 1987.6666
$ws.Range("I102").Value = 1305.32
$ws.Range("K102").Value = 1305.32
$ws.Range("M102").Value = 316.6800000000001
$ws.Range("H132").Value = 2348.3
$ws.Range("I132").Value = 2212
$ws.Range("K132").Value = 6636
$ws.Range("M132").Value = -4106

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4721.1333
$ws.Range("I7").Value = 3901.0908
$ws.Range("K7").Value = 3901.0908
$ws.Range("M7").Value = -3789.0908
$ws.Range("H40").Value = 15919
$ws.Range("I40").Value = 19003.5
$ws.Range("J40").Value = 6665.5
$ws.Range("K40").Value = 19003.5
$ws.Range("L40").Value = 6665.5
$ws.Range("M40").Value = -18867.5
$ws.Range("N40").Value = -6937.5
$ws.Range("H126").Value = 4721.1333
$ws.Range("I126").Value = 3901.0908
$ws.Range("K126").Value = 11703.2724
$ws.Range("M126").Value = -9233.2724
$ws.Range("H132").Value = 6174.5415
$ws.Range("I132").Value = 3752.4285
$ws.Range("K132").Value = 11257.2855
$ws.Range("M132").Value = -8727.2855
$ws.Range("H136").Value = 2602.5469
$ws.Range("I136").Value = 2329.1018
$ws.Range("J136").Value = 5829.2
$ws.Range("K136").Value = 6987.305399999999
$ws.Range("L136").Value = 17487.6
$ws.Range("M136").Value = -4437.305399999999
$ws.Range("N136").Value = -22587.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18523852
$ws.Range("I81").Value = 5332.3335
$ws.Range("K81").Value = 10664.667
$ws.Range("M81").Value = -9603.666999999999
$ws.Range("H84").Value = 18523852
$ws.Range("I84").Value = 5332.3335
$ws.Range("K84").Value = 53323.335
$ws.Range("M84").Value = -48019.335
$ws.Range("H122").Value = 269906.4
$ws.Range("I122").Value = 2353.8125
$ws.Range("J122").Value = 1126074.8
$ws.Range("K122").Value = 7061.4375
$ws.Range("L122").Value = 3378224.4
$ws.Range("M122").Value = -4611.4375
$ws.Range("N122").Value = -3383124.4
$ws.Range("H132").Value = 3121.6667
$ws.Range("I132").Value = 2156.5
$ws.Range("K132").Value = 6469.5
$ws.Range("M132").Value = -3939.5
